# Update pl_mw.xlsx results for case with 380 kV: refresh the line power-flow
# results table (rows 2-25, columns B,C,D,E,G,H,J,K) with newly computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7303842732983412
$ws.Range("C2").Value = 0.06990963580585685
$ws.Range("D2").Value = 0.6767035258003489
$ws.Range("E2").Value = 0.2755673430226437
$ws.Range("G2").Value = 3.060189984765742
$ws.Range("H2").Value = 2.307321076619587
$ws.Range("J2").Value = 0.1430532246694796
$ws.Range("K2").Value = 0.7220761204327459
$ws.Range("B3").Value = 0.7038350451051656
$ws.Range("C3").Value = 0.06679977807573323
$ws.Range("D3").Value = 0.6674835966361456
$ws.Range("E3").Value = 0.2709806420593637
$ws.Range("G3").Value = 2.988258369005621
$ws.Range("H3").Value = 2.27751785654462
$ws.Range("J3").Value = 0.139968016515013
$ws.Range("K3").Value = 0.6948936435151438
$ws.Range("B4").Value = 0.6881313307902417
$ws.Range("C4").Value = 0.06497022058290725
$ws.Range("D4").Value = 0.6621814344264862
$ws.Range("E4").Value = 0.2683175004330494
$ws.Range("G4").Value = 2.944937172142261
$ws.Range("H4").Value = 2.259796618350009
$ws.Range("J4").Value = 0.1381581173780191
$ws.Range("K4").Value = 0.6788364726088219
$ws.Range("B5").Value = 0.6818819778537488
$ws.Range("C5").Value = 0.06424466735195722
$ws.Range("D5").Value = 0.6601108838342782
$ws.Range("E5").Value = 0.2672706678702923
$ws.Range("G5").Value = 2.927494433010708
$ws.Range("H5").Value = 2.25271983344112
$ws.Range("J5").Value = 0.1374417309590399
$ws.Range("K5").Value = 0.67245188780376
$ws.Range("B6").Value = 0.6808533327882458
$ws.Range("C6").Value = 0.06412539561493702
$ws.Range("D6").Value = 0.6597725101003675
$ws.Range("E6").Value = 0.2670991602055466
$ws.Range("G6").Value = 2.924610779489797
$ws.Range("H6").Value = 2.251553464802981
$ws.Range("J6").Value = 0.1373240517018459
$ws.Range("K6").Value = 0.6714013168384838
$ws.Range("B7").Value = 0.6880464426429569
$ws.Range("C7").Value = 0.06496035462828331
$ws.Range("D7").Value = 0.6621531455143099
$ws.Range("E7").Value = 0.2683032270306214
$ws.Range("G7").Value = 2.944701080932845
$ws.Range("H7").Value = 2.259700593020426
$ws.Range("J7").Value = 0.1381483703447373
$ws.Range("K7").Value = 0.6787497251424099
$ws.Range("B8").Value = 0.7211059062411493
$ws.Range("C8").Value = 0.06882071772557197
$ws.Range("D8").Value = 0.6734499194764396
$ws.Range("E8").Value = 0.2739540183715192
$ws.Range("G8").Value = 3.035211627523353
$ws.Range("H8").Value = 2.296924579832393
$ws.Range("J8").Value = 0.1419718782165589
$ws.Range("K8").Value = 0.7125720169727288
$ws.Range("B9").Value = 0.7906931984257426
$ws.Range("C9").Value = 0.0770294726396088
$ws.Range("D9").Value = 0.6984586076710286
$ws.Range("E9").Value = 0.2862547639414785
$ws.Range("G9").Value = 3.219486578570638
$ws.Range("H9").Value = 2.374539148810129
$ws.Range("J9").Value = 0.1501433988455076
$ws.Range("K9").Value = 0.7839398942285243
$ws.Range("B10").Value = 0.8447490603749372
$ws.Range("C10").Value = 0.08345698068525564
$ws.Range("D10").Value = 0.7185876003397595
$ws.Range("E10").Value = 0.2960435428454389
$ws.Range("G10").Value = 3.359132915113491
$ws.Range("H10").Value = 2.434429349751127
$ws.Range("J10").Value = 0.1565636993857566
$ws.Range("K10").Value = 0.8394834505520237
$ws.Range("B11").Value = 0.869983740137485
$ws.Range("C11").Value = 0.08646877457145763
$ws.Range("D11").Value = 0.7281292354047935
$ws.Range("E11").Value = 0.3006617131284699
$ws.Range("G11").Value = 3.423614286509348
$ws.Range("H11").Value = 2.462309014251503
$ws.Range("J11").Value = 0.1595763122579257
$ws.Range("K11").Value = 0.8654352204852387
$ws.Range("B12").Value = 0.8796325633464335
$ws.Range("C12").Value = 0.08762202156832188
$ws.Range("D12").Value = 0.7317979643982255
$ws.Range("E12").Value = 0.3024343754141867
$ws.Range("G12").Value = 3.448171050664996
$ws.Range("H12").Value = 2.472958472071866
$ws.Range("J12").Value = 0.1607304342860232
$ws.Range("K12").Value = 0.8753614930674303
$ws.Range("B13").Value = 0.8775503750730138
$ws.Range("C13").Value = 0.08737308056764448
$ws.Range("D13").Value = 0.7310053658710274
$ws.Range("E13").Value = 0.3020515379304101
$ws.Range("G13").Value = 3.44287610042096
$ws.Range("H13").Value = 2.470660816684585
$ws.Range("J13").Value = 0.1604812801049462
$ws.Range("K13").Value = 0.8732192869077835
$ws.Range("B14").Value = 0.8707756888583162
$ws.Range("C14").Value = 0.08656339682069358
$ws.Range("D14").Value = 0.7284299502308045
$ws.Range("E14").Value = 0.3008070722918816
$ws.Range("G14").Value = 3.42563178957505
$ws.Range("H14").Value = 2.463183303155404
$ws.Range("J14").Value = 0.1596709953172848
$ws.Range("K14").Value = 0.8662498762310236
$ws.Range("B15").Value = 0.8666381191149242
$ws.Range("C15").Value = 0.0860691052985203
$ws.Range("D15").Value = 0.7268596693382108
$ws.Range("E15").Value = 0.3000479117798065
$ws.Range("G15").Value = 3.415087311597404
$ws.Range("H15").Value = 2.458615115962857
$ws.Range("J15").Value = 0.1591764084672604
$ws.Range("K15").Value = 0.8619938039631734
$ws.Range("B16").Value = 0.8431129098889301
$ws.Range("C16").Value = 0.0832619305643334
$ws.Range("D16").Value = 0.7179717913383001
$ws.Range("E16").Value = 0.2957450678382543
$ws.Range("G16").Value = 3.354938275656139
$ws.Range("H16").Value = 2.432620189891765
$ws.Range("J16").Value = 0.1563686761382712
$ws.Range("K16").Value = 0.8378012568861379
$ws.Range("B17").Value = 0.8288462351881378
$ws.Range("C17").Value = 0.08156240595475595
$ws.Range("D17").Value = 0.7126180610341351
$ws.Range("E17").Value = 0.293147799512063
$ws.Range("G17").Value = 3.318284707459753
$ws.Range("H17").Value = 2.416836342484117
$ws.Range("J17").Value = 0.1546698419674755
$ws.Range("K17").Value = 0.8231355946256826
$ws.Range("B18").Value = 0.8207010541698594
$ws.Range("C18").Value = 0.08059315020517488
$ws.Range("D18").Value = 0.7095749504439652
$ws.Range("E18").Value = 0.2916694629677536
$ws.Range("G18").Value = 3.29729238544553
$ws.Range("H18").Value = 2.40781764855285
$ws.Range("J18").Value = 0.1537013705509054
$ws.Range("K18").Value = 0.8147647163189617
$ws.Range("B19").Value = 0.8179536396212939
$ws.Range("C19").Value = 0.08026639328690521
$ws.Range("D19").Value = 0.7085508198770469
$ws.Range("E19").Value = 0.2911715900674352
$ws.Range("G19").Value = 3.290200125801675
$ws.Range("H19").Value = 2.404774316667726
$ws.Range("J19").Value = 0.1533749464261831
$ws.Range("K19").Value = 0.8119415323284045
$ws.Range("B20").Value = 0.8303586694167677
$ws.Range("C20").Value = 0.08174246699401522
$ws.Range("D20").Value = 0.7131842255411414
$ws.Range("E20").Value = 0.2934226737121577
$ws.Range("G20").Value = 3.322177232201625
$ws.Range("H20").Value = 2.418510370461661
$ws.Range("J20").Value = 0.1548497896049952
$ws.Range("K20").Value = 0.8246901087884453
$ws.Range("B21").Value = 0.8727630517631155
$ws.Range("C21").Value = 0.08680087376842494
$ws.Range("D21").Value = 0.7291849040281591
$ws.Range("E21").Value = 0.3011719535739061
$ws.Range("G21").Value = 3.430693074096041
$ws.Range("H21").Value = 2.465377126927649
$ws.Range("J21").Value = 0.1599086336735667
$ws.Range("K21").Value = 0.8682942733944685
$ws.Range("B22").Value = 0.901018992787499
$ws.Range("C22").Value = 0.09018118043573509
$ws.Range("D22").Value = 0.739965965095422
$ws.Range("E22").Value = 0.3063756892144269
$ws.Range("G22").Value = 3.502425827426009
$ws.Range("H22").Value = 2.496544025495098
$ws.Range("J22").Value = 0.1632925068199569
$ws.Range("K22").Value = 0.8973687973852691
$ws.Range("B23").Value = 0.8858885377771912
$ws.Range("C23").Value = 0.08837020824341835
$ws.Range("D23").Value = 0.7341822320054519
$ws.Range("E23").Value = 0.3035855914712684
$ws.Range("G23").Value = 3.464065913191348
$ws.Range("H23").Value = 2.47986032554428
$ws.Range("J23").Value = 0.1614793392423195
$ws.Range("K23").Value = 0.8817982618600979
$ws.Range("B24").Value = 0.8296747206286739
$ws.Range("C24").Value = 0.0816610370377191
$ws.Range("D24").Value = 0.7129281541235173
$ws.Range("E24").Value = 0.2932983567855274
$ws.Range("G24").Value = 3.320417171757128
$ws.Range("H24").Value = 2.417753369104389
$ws.Range("J24").Value = 0.1547684097096607
$ws.Range("K24").Value = 0.8239871241887613
$ws.Range("B25").Value = 0.7713556242904076
$ws.Range("C25").Value = 0.07473968252541852
$ws.Range("D25").Value = 0.6913859221455709
$ws.Range("E25").Value = 0.2827957117605564
$ws.Range("G25").Value = 3.16889629469361
$ws.Range("H25").Value = 2.353042643237757
$ws.Range("J25").Value = 0.1478600522970979
$ws.Range("K25").Value = 0.764089539129543
